$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F3").Value = 2157
    $ws.Range("F4").Value = 313
    $ws.Range("F5").Value = 74
    $ws.Range("F6").Value = 6404
    $ws.Range("F7").Value = 275
}
